$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Taurean Prince -> Stephen Curry
$ws.Range("A5").Value = "Stephen Curry"
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "Golden State Warriors"

# Row 12: Keegan Murray -> Brandon Clarke
$ws.Range("A12").Value = "Brandon Clarke"
$ws.Range("B12").Value = "PF,C"
$ws.Range("C12").Value = "Memphis Grizzlies"

# Row 13: Stephen Curry -> Keegan Murray
$ws.Range("A13").Value = "Keegan Murray"
$ws.Range("B13").Value = "SF,PF"
$ws.Range("C13").Value = "Sacramento Kings"
